$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells we are about to rewrite to remain plain text,
# matching the source workbook where these are inline strings (not numbers).
$priceCells = @(
    "D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11",
    "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21",
    "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D31", "D32",
    "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43",
    "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values row by row
# Row 2
$ws.Range("D2").Value = "29.397.28"
$ws.Range("E2").Value = "  -1.90%  "
# Row 3
$ws.Range("D3").Value = "1.997.52"
$ws.Range("E3").Value = "  -5.52%  "
# Row 4
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  -0.02%  "
# Row 5
$ws.Range("D5").Value = "331.31"
$ws.Range("E5").Value = "  -4.34%  "
# Row 6
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  -0.08%  "
# Row 7
$ws.Range("D7").Value = "0.4941"
$ws.Range("E7").Value = "  -4.72%  "
# Row 8
$ws.Range("D8").Value = "0.4183"
$ws.Range("E8").Value = "  -5.86%  "
# Row 9
$ws.Range("D9").Value = "53.37"
$ws.Range("E9").Value = "  -0.50%  "
# Row 10
$ws.Range("D10").Value = "0.08782"
$ws.Range("E10").Value = "  -6.32%  "
# Row 11
$ws.Range("D11").Value = "1.112"
$ws.Range("E11").Value = "  -5.93%  "
# Row 12
$ws.Range("D12").Value = "2.201.93"
$ws.Range("E12").Value = "  +4.29%  "
# Row 13
$ws.Range("D13").Value = "23.09"
$ws.Range("E13").Value = "  -8.37%  "
# Row 14
$ws.Range("D14").Value = "8.077"
$ws.Range("E14").Value = "  -5.25%  "
# Row 15
$ws.Range("D15").Value = "6.471"
$ws.Range("E15").Value = "  -6.40%  "
# Row 16
$ws.Range("D16").Value = "96.11"
$ws.Range("E16").Value = "  -6.63%  "
# Row 17
$ws.Range("D17").Value = "1.009"
$ws.Range("E17").Value = "  +0.09%  "
# Row 18
$ws.Range("D18").Value = "0.00001104"
$ws.Range("E18").Value = "  -5.08%  "
# Row 19
$ws.Range("D19").Value = "0.06636"
$ws.Range("E19").Value = "  -0.88%  "
# Row 20
$ws.Range("D20").Value = "19.46"
$ws.Range("E20").Value = "  -9.40%  "
# Row 21
$ws.Range("D21").Value = "1.008"
$ws.Range("E21").Value = "  -0.01%  "
# Row 22
$ws.Range("D22").Value = "5.956"
$ws.Range("E22").Value = "  -5.40%  "
# Row 23
$ws.Range("D23").Value = "29.462.60"
$ws.Range("E23").Value = "  -1.79%  "
# Row 24
$ws.Range("D24").Value = "11.76"
$ws.Range("E24").Value = "  -7.49%  "
# Row 25
$ws.Range("D25").Value = "2.284"
$ws.Range("E25").Value = "  -1.48%  "
# Row 26
$ws.Range("D26").Value = "2.354.70"
$ws.Range("E26").Value = "  -0.19%  "
# Row 27
$ws.Range("D27").Value = "6.672"
$ws.Range("E27").Value = "  +0.86%  "
# Row 28
$ws.Range("D28").Value = "157.56"
$ws.Range("E28").Value = "  -3.02%  "
# Row 29
$ws.Range("D29").Value = "20.51"
$ws.Range("E29").Value = "  -7.16%  "
# Row 30
$ws.Range("E30").Value = "  -7.59%  "
# Row 31
$ws.Range("D31").Value = "126.69"
$ws.Range("E31").Value = "  -5.56%  "
# Row 32
$ws.Range("D32").Value = "1.046"
$ws.Range("E32").Value = "  -8.96%  "
# Row 33
$ws.Range("D33").Value = "0.09912"
$ws.Range("E33").Value = "  -6.10%  "
# Row 34
$ws.Range("D34").Value = "1.551"
$ws.Range("E34").Value = "  -13.05%  "
# Row 35
$ws.Range("D35").Value = "5.805"
$ws.Range("E35").Value = "  -7.04%  "
# Row 36
$ws.Range("D36").Value = "3.782"
$ws.Range("E36").Value = "  -4.78%  "
# Row 37
$ws.Range("D37").Value = "9.606"
$ws.Range("E37").Value = "  -10.99%  "
# Row 38
$ws.Range("D38").Value = "0.02443"
$ws.Range("E38").Value = "  -6.40%  "
# Row 39
$ws.Range("D39").Value = "0.06364"
$ws.Range("E39").Value = "  -7.17%  "
# Row 40
$ws.Range("E40").Value = "  -3.82%  "
# Row 41
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.6481"
$ws.Range("E41").Value = "  -8.63%  "
# Row 42
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "11.73"
$ws.Range("E42").Value = "  -7.62%  "
# Row 43
$ws.Range("D43").Value = "0.2060"
$ws.Range("E43").Value = "  -7.96%  "
# Row 44
$ws.Range("D44").Value = "1.008"
$ws.Range("E44").Value = "  +0.07%  "
# Row 45
$ws.Range("D45").Value = "0.6302"
$ws.Range("E45").Value = "  -7.82%  "
# Row 46
$ws.Range("D46").Value = "2.194"
$ws.Range("E46").Value = "  -7.22%  "
# Row 47
$ws.Range("D47").Value = "13.37"
$ws.Range("E47").Value = "  -8.55%  "
# Row 48
$ws.Range("D48").Value = "1.255"
$ws.Range("E48").Value = "  -1.45%  "
# Row 49
$ws.Range("D49").Value = "3.557"
$ws.Range("E49").Value = "  -2.04%  "
# Row 50
$ws.Range("D50").Value = "0.06986"
$ws.Range("E50").Value = "  -1.73%  "
# Row 51
$ws.Range("D51").Value = "1.147"
$ws.Range("E51").Value = "  -3.07%  "

# Restore default (Normal) style on the price cells so no stray number format remains
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
